$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "IP" column (H) formatting into the two new columns (I, J)
# so the new header cells (I1, J1) pick up the same bold/bordered/centered
# style as the other header cells, and the data cells (I2:I5, J2:J5) stay
# in the default (unstyled) format, matching column H.
$ws.Range("H1:H5").Copy()
$ws.Range("I1:I5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1:H5").Copy()
$ws.Range("J1:J5").PasteSpecial(-4122)  # xlPasteFormats

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 5
